# FacebookReg.xlsx edit:
#  - RegTestData: add "status"/"pass" columns (E,F) and a new data row (satyam)
#  - HomePage: populate sample "newdata"/"hello" cells
#  - add a new "satyam" worksheet at the end with the same sample cells

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # RegTestData
$ws2 = $wb.Worksheets.Item(2)   # HomePage

# --- RegTestData: new "status" columns for every existing row ---
$ws1.Range("E1").Value = "status"
$ws1.Range("F1").Value = "status"

$ws1.Range("E2").Value = "pass"
$ws1.Range("F2").Value = "pass"

$ws1.Range("E3").Value = "pass"
$ws1.Range("F3").Value = "pass"

$ws1.Range("E4").Value = "pass"
$ws1.Range("F4").Value = "pass"

# --- RegTestData: new data row for "satyam" ---
$ws1.Range("A5").Value = "satyam"
$ws1.Range("B5").Value = "s"
$ws1.Range("C5").Value = 3425252625
$ws1.Range("D5").Value = "Q@123"
$ws1.Range("E5").Value = "pass"
$ws1.Range("F5").Value = "pass"

# size the two new columns to fit their (short) content
$ws1.Columns.Item(5).ColumnWidth = 5
$ws1.Columns.Item(6).ColumnWidth = 5

# --- HomePage: sample data used by the new DDT code ---
$ws2.Range("A1").Value = "newdata"
$ws2.Range("B1").Value = "newdata"
$ws2.Range("A2").Value = "hello"
$ws2.Range("B3").Value = "hello"
$ws2.Columns.Item(1).ColumnWidth = 7.9
$ws2.Columns.Item(2).ColumnWidth = 7.9
[void]$ws2.Range("P8").Select()

# --- add the new "satyam" worksheet at the end of the workbook ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$newSheet.Name = "satyam"
$newSheet.Range("A1").Value = "newdata"
$newSheet.Range("A3").Value = "hello"
$newSheet.Columns.Item(1).ColumnWidth = 7.9

# leave RegTestData as the active sheet/selection, like the source file
[void]$ws1.Select()
[void]$ws1.Range("A1").Select()
